# Site updated: 2022-05-17 16:37:31
#
# Append three new article rows (44-46) to the tracking sheet, mirroring
# the existing row layout:
#   A: publish id (unix timestamp, number)
#   B: publish date (plain text, e.g. "2022-05-17" - NOT a date value)
#   C: article title
#   D: source account
#   E: cover image url
#   F: article link

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 44; A = 1652766746; B = "2022-05-17"; C = "【招聘】Sophoton | 识光芯科在招岗位"; D = "JIcareer"; E = "https://mmbiz.qlogo.cn/mmbiz_jpg/PB5L45ia5moSxJYWoqYh7Mcrp4pecJNeJiamkvhgNySd3PIVbsKgPZxBXo2tjcB9FibIF1XMvb6FoogQvdgCgDo2A/0?wx_fmt=jpeg"; F = "http://mp.weixin.qq.com/s?__biz=MzI5MTcwMDg2Mg==&mid=2247491691&idx=1&sn=70aad1997dc0ba16c809ab3d1e53b883&chksm=ec0e097adb79806cd753c4950d1309b9ad5f8e3c0a65d1a08e98030376caab819132b121ff3c#rd" },
    @{ Row = 45; A = 1652708097; B = "2022-05-16"; C = "【活动回顾】GDP×保研workshop圆满结束！"; D = "JIcareer"; E = "https://mmbiz.qlogo.cn/mmbiz_jpg/PB5L45ia5moSictu4icWQJwfw0Myibnt9xKLjlv13GbzAt9D2Vp7My36zeOH1cOQAYD4l68OZG638slxbzWuh4liczA/0?wx_fmt=jpeg"; F = "http://mp.weixin.qq.com/s?__biz=MzI5MTcwMDg2Mg==&mid=2247491681&idx=1&sn=9d27cb02f3aacbe6bfd5e8a4a3540e6a&chksm=ec0e0970db798066032ba2bb809e8c29473206ce54164131221d23bbc2a7c1c8907e8f8ab307#rd" },
    @{ Row = 46; A = 1652711615; B = "2022-05-16"; C = "青年大学习：学习习近平总书记在庆祝中国共产主义青年团成立100周年大会上的重要讲话精神"; D = "JI青团"; E = "https://mmbiz.qlogo.cn/mmbiz_jpg/QfDapvG9u4Dib2WjxQn9Bgn5VticgfEfdKibYtX5BQA5zbt9gemyzTrFYp1oOzj1V6MyBUjqsV7mhqrKoCntc8ByQ/0?wx_fmt=jpeg"; F = "http://mp.weixin.qq.com/s?__biz=MzUyMzMyNTY0OQ==&mid=2247486248&idx=1&sn=92d251a227a6a0190374ff5a283944af&chksm=fa3f1296cd489b8003c01dfc403cce50dfd33c649544333390b6e15f72e1e3b738b7ef524719#rd" }
)

$lastRow = 43
foreach ($r in $rows) {
    $row = $r.Row

    # Carry the row's cell formatting down from the row above it, so the
    # new rows look like a natural continuation of the table.
    $ws.Range("A" + $lastRow + ":F" + $lastRow).Copy()
    $ws.Range("A" + $row + ":F" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.A

    # Column B holds a date-looking string ("2022-05-17") that must stay
    # plain text like the rest of the column, instead of being
    # auto-parsed into a date serial number. Force text entry via
    # NumberFormat, assign the value, then re-apply the row's normal
    # formatting (font/number format) on top - this does not disturb the
    # already-committed text value/type.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Range("B" + $lastRow).Copy()
    $ws.Range("B" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F

    $lastRow = $row
}

# Extend the sheet's active selection to cover the newly added rows,
# matching the widened used range (A1:F46).
$null = $ws.Range("A1:F46").Select()
